# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update "last updated" timestamp text (A1, uses shared string) ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 14:22"

# --- Swap order of "Corea del Sur" and "Suecia" rows (21-24 block) ---
# Before: row21=Israel, row22=Corea del Sur, row23=Suecia, row24=Irlanda
# After:  row21=Israel, row22=Suecia,        row23=Corea del Sur, row24=Irlanda
# Row22 gets freshly updated data, row23 gets the data that used to belong
# to the old "Corea del Sur" row (i.e. the old row22 values).
$ws.Range("A22").Value = "Suecia"
$ws.Range("B22").Value = 10948
$ws.Range("C22").Value = 465
$ws.Range("D22").Value = 381
$ws.Range("E22").Value = 9648
$ws.Range("F22").Value = 859
$ws.Range("G22").Value = 20
$ws.Range("H22").Value = 919

$ws.Range("A23").Value = "Corea del Sur"
$ws.Range("B23").Value = 10537
$ws.Range("C23").Value = 25
$ws.Range("D23").Value = 7447
$ws.Range("E23").Value = 2873
$ws.Range("F23").Value = 55
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 217

# --- Row 14 data updates ---
$ws.Range("B14").Value = 26551
$ws.Range("C14").Value = 964
$ws.Range("E14").Value = 23478
$ws.Range("G14").Value = 86
$ws.Range("H14").Value = 2823

# --- Row 17 data updates ---
$ws.Range("E17").Value = 20904
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = 1241

# --- Row 20 data updates ---
$ws.Range("B20").Value = 13998
$ws.Range("C20").Value = 53
$ws.Range("E20").Value = 6287

# --- Row 34 data updates ---
$ws.Range("B34").Value = 6318
$ws.Range("C34").Value = 144
$ws.Range("D34").Value = 2235
$ws.Range("E34").Value = 3798
$ws.Range("F34").Value = 100
$ws.Range("G34").Value = 12
$ws.Range("H34").Value = 285

# --- Row 62 data updates ---
$ws.Range("B62").Value = 1650
$ws.Range("C62").Value = 50
$ws.Range("D62").Value = 400
$ws.Range("E62").Value = 1225
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 25

# --- Row 83 data updates ---
$ws.Range("F83").Value = 89

# --- Row 87 data updates ---
$ws.Range("B87").Value = 655
$ws.Range("C87").Value = 4
$ws.Range("E87").Value = 634
